# Natmi following Dr Hou advice:
# Re-run the Bmp2-Acvr2a ligand-receptor pair analysis with "ECs" added as an
# additional sending/target cluster. This rewrites the data rows (A2:T7) with
# updated statistics for the FAPs/sCs combinations and appends three new rows
# (A8:T10) for the sCs-sending-cluster combinations that complete the
# ECs/FAPs/sCs x ECs/FAPs/sCs grid.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Bmp2","Acvr2a","ECs",2,0.6666666666666666,9.163165666666666,27.489497,0.2800251397703982,0.2800251397703982,3,1,7.940402333333334,23.821207,0.1931648990487216,0.1931648990487216,72.75922204031988,654.832998362879,0.05409102785485313,0.05409102785485314),
    @("ECs","Bmp2","Acvr2a","FAPs",2,0.6666666666666666,9.163165666666666,27.489497,0.2800251397703982,0.2800251397703982,3,1,26.95384733333333,80.861542,0.6557019380820612,0.6557019380820612,246.9825684693749,2222.843116224374,0.1836130268591502,0.1836130268591502),
    @("ECs","Bmp2","Acvr2a","sCs",2,0.6666666666666666,9.163165666666666,27.489497,0.2800251397703982,0.2800251397703982,3,1,6.212609666666666,18.637829,0.1511331628692172,0.1511331628692172,56.92717159800144,512.344544382013,0.04232108505639488,0.0423210850563949),
    @("FAPs","Bmp2","Acvr2a","ECs",3,1,12.06704733333333,36.201142,0.3687673822623249,0.3687673822623249,3,1,7.940402333333334,23.821207,0.1931648990487216,0.1931648990487216,95.81721080204379,862.354897218394,0.07123291416716332,0.07123291416716332),
    @("FAPs","Bmp2","Acvr2a","FAPs",3,1,12.06704733333333,36.201142,0.3687673822623249,0.3687673822623249,3,1,26.95384733333333,80.861542,0.6557019380820612,0.6557019380820612,325.2533515867738,2927.280164280964,0.2418014872508547,0.2418014872508547),
    @("FAPs","Bmp2","Acvr2a","sCs",3,1,12.06704733333333,36.201142,0.3687673822623249,0.3687673822623249,3,1,6.212609666666666,18.637829,0.1511331628692172,0.1511331628692172,74.96785491119088,674.710694200718,0.05573298084430681,0.05573298084430681),
    @("sCs","Bmp2","Acvr2a","ECs",3,1,11.49244066666667,34.477322,0.3512074779672769,0.351207477967277,3,1,7.940402333333334,23.821207,0.1931648990487216,0.1931648990487216,91.2546026852949,821.2914241676541,0.06784095702670517,0.06784095702670519),
    @("sCs","Bmp2","Acvr2a","FAPs",3,1,11.49244066666667,34.477322,0.3512074779672769,0.351207477967277,3,1,26.95384733333333,80.861542,0.6557019380820612,0.6557019380820612,309.7654912167249,2787.889420950524,0.2302874239720563,0.2302874239720563),
    @("sCs","Bmp2","Acvr2a","sCs",3,1,11.49244066666667,34.477322,0.3512074779672769,0.351207477967277,3,1,6.212609666666666,18.637829,0.1511331628692172,0.1511331628692172,71.39804797932645,642.582431813938,0.05307909696851546,0.05307909696851548)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}
